$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Writing the literal word True/False (or other auto-typed tokens) via
    # .Value always gets coerced to a Boolean by the engine; routing the
    # text through a formula and then collapsing it to a static value via
    # copy / paste-special (values only) keeps it as a genuine shared string.
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# ----- header row -----
$ws.Range("A1").Value = "CHARACTER_CLASS"
$ws.Range("B1").Value = "FRIEND"
$ws.Range("C1").Value = "LIFE"
$ws.Range("D1").Value = "DISTANCE"
$ws.Range("E1").Value = "FOLLOW"

# ----- row 2 : ORC / not-friend / close life band -----
$ws.Range("A2").Value = "ORC"
Set-TextValue $ws.Range("B2") "False"
$ws.Range("C2").Value = ">=50"
$ws.Range("D2").Value = "<=60"
Set-TextValue $ws.Range("E2") "True"

# ----- row 3 : ORC / friend / low life band -----
$ws.Range("A3").Value = "ORC"
Set-TextValue $ws.Range("B3") "True"
$ws.Range("C3").Value = "<=0"
$ws.Range("D3").Value = ">=60"
Set-TextValue $ws.Range("E3") "False"

# ----- row 4 : MONKEY / not-friend / close life band -----
$ws.Range("A4").Value = "MONKEY"
Set-TextValue $ws.Range("B4") "False"
$ws.Range("C4").Value = ">=50"
$ws.Range("D4").Value = "<=60"
Set-TextValue $ws.Range("E4") "True"

# ----- row 5 : MONKEY / friend / low life band -----
$ws.Range("A5").Value = "MONKEY"
Set-TextValue $ws.Range("B5") "True"
$ws.Range("C5").Value = "<=0"
$ws.Range("D5").Value = ">=60"
Set-TextValue $ws.Range("E5") "False"

# ----- bring the new column E cells onto the same look as the rest -----
# (right aligned Arial, matching style index already used by columns A-D)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# E2 gets its own (slightly) distinct style/font entry
$donor = $ws.Range("Z100")
$donor.Value = "x"
$ws.Range("D1").Copy()
$donor.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$donor.Font.Bold = $true
$donor.Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$donor.Clear()

# ----- column widths -----
$ws.Columns.Item(1).ColumnWidth = 19.65
$ws.Columns.Item(2).ColumnWidth = 28.8

# ----- selection / view -----
$ws.Range("A3").Select()

# ----- page setup -----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
